$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the name in row 2
$ws.Range("B2").Value = "Ahmed Al Mansoori"

# Update employment status for rows 3 and 4 (swapped)
$ws.Range("D3").Value = "Self-employed"
$ws.Range("D4").Value = "Unemployed"

# Update Family Members counts
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 0

# Move the active selection to E8
$ws.Range("E8").Select()
